# Investigation.docx edit:
#   "enjoy watching arts"  ->  "Enjoy watching arts"
# (capitalize the leading "e" of "enjoy", matching the capitalization style
#  already used elsewhere in the document, e.g. "Choose ..." / "Each way...")
#
# Word marks the location of the most recent edit with a hidden "_GoBack"
# bookmark. Making this small, single-character edit moves that bookmark
# from its previous location to right after the newly-typed "E".

$d = $word.ActiveDocument

# Locate the word to fix by scanning paragraphs (avoids Find.Execute so we
# don't disturb unrelated document statistics).
$target = "enjoy watching arts"
$startPos = -1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    $idx = $t.IndexOf($target)
    if ($idx -ge 0) {
        $startPos = $p.Range.Start + $idx
        break
    }
}

if ($startPos -ge 0) {
    # Remove the stale _GoBack bookmark (Word relocates this to the newest edit point).
    try {
        $old = $d.Bookmarks("_GoBack")
        $old.Delete()
    } catch {
        # no existing _GoBack bookmark - nothing to remove
    }

    # Replace just the leading "e" with "E", leaving the rest of the word/run untouched.
    $firstChar = $d.Range($startPos, $startPos + 1)
    $firstChar.Text = "E"

    # Re-create the _GoBack bookmark at the new edit point, right after the "E".
    $editPoint = $d.Range($startPos + 1, $startPos + 1)
    $d.Bookmarks.Add("_GoBack", $editPoint)
}
